$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Sema4g"
$ws.Cells.Item(2,3).Value = "Plxnb2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2.0
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.4290743333333333
$ws.Cells.Item(2,8).Value = 1.287223
$ws.Cells.Item(2,9).Value = 0.1093614795344676
$ws.Cells.Item(2,10).Value = 0.1093614795344676
$ws.Cells.Item(2,11).Value = 3.0
$ws.Cells.Item(2,12).Value = 1.0
$ws.Cells.Item(2,13).Value = 8.806900666666666
$ws.Cells.Item(2,14).Value = 26.420702
$ws.Cells.Item(2,15).Value = 0.1733678197953833
$ws.Cells.Item(2,16).Value = 0.1733678197953834
$ws.Cells.Item(2,17).Value = 3.778815032282889
$ws.Cells.Item(2,18).Value = 34.009335290546
$ws.Cells.Item(2,19).Value = 0.01895976127648808
$ws.Cells.Item(2,20).Value = 0.01895976127648808

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Sema4g"
$ws.Cells.Item(3,3).Value = "Plxnb2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2.0
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.4290743333333333
$ws.Cells.Item(3,8).Value = 1.287223
$ws.Cells.Item(3,9).Value = 0.1093614795344676
$ws.Cells.Item(3,10).Value = 0.1093614795344676
$ws.Cells.Item(3,11).Value = 3.0
$ws.Cells.Item(3,12).Value = 1.0
$ws.Cells.Item(3,13).Value = 18.76689066666667
$ws.Cells.Item(3,14).Value = 56.30067200000001
$ws.Cells.Item(3,15).Value = 0.3694347242421866
$ws.Cells.Item(3,16).Value = 0.3694347242421866
$ws.Cells.Item(3,17).Value = 8.052391101539557
$ws.Cells.Item(3,18).Value = 72.47151991385601
$ws.Cells.Item(3,19).Value = 0.04040192803453356
$ws.Cells.Item(3,20).Value = 0.04040192803453356

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Sema4g"
$ws.Cells.Item(4,3).Value = "Plxnb2"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2.0
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.4290743333333333
$ws.Cells.Item(4,8).Value = 1.287223
$ws.Cells.Item(4,9).Value = 0.1093614795344676
$ws.Cells.Item(4,10).Value = 0.1093614795344676
$ws.Cells.Item(4,11).Value = 3.0
$ws.Cells.Item(4,12).Value = 1.0
$ws.Cells.Item(4,13).Value = 23.225144
$ws.Cells.Item(4,14).Value = 69.675432
$ws.Cells.Item(4,15).Value = 0.4571974559624301
$ws.Cells.Item(4,16).Value = 0.4571974559624301
$ws.Cells.Item(4,17).Value = 9.965313178370668
$ws.Cells.Item(4,18).Value = 89.687818605336
$ws.Cells.Item(4,19).Value = 0.04999979022344594
$ws.Cells.Item(4,20).Value = 0.04999979022344594

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Sema4g"
$ws.Cells.Item(5,3).Value = "Plxnb2"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3.0
$ws.Cells.Item(5,6).Value = 1.0
$ws.Cells.Item(5,7).Value = 1.525446
$ws.Cells.Item(5,8).Value = 4.576338
$ws.Cells.Item(5,9).Value = 0.3888021691111845
$ws.Cells.Item(5,10).Value = 0.3888021691111845
$ws.Cells.Item(5,11).Value = 3.0
$ws.Cells.Item(5,12).Value = 1.0
$ws.Cells.Item(5,13).Value = 8.806900666666666
$ws.Cells.Item(5,14).Value = 26.420702
$ws.Cells.Item(5,15).Value = 0.1733678197953833
$ws.Cells.Item(5,16).Value = 0.1733678197953834
$ws.Cells.Item(5,17).Value = 13.434451394364
$ws.Cells.Item(5,18).Value = 120.910062549276
$ws.Cells.Item(5,19).Value = 0.067405784390522
$ws.Cells.Item(5,20).Value = 0.06740578439052201

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Sema4g"
$ws.Cells.Item(6,3).Value = "Plxnb2"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3.0
$ws.Cells.Item(6,6).Value = 1.0
$ws.Cells.Item(6,7).Value = 1.525446
$ws.Cells.Item(6,8).Value = 4.576338
$ws.Cells.Item(6,9).Value = 0.3888021691111845
$ws.Cells.Item(6,10).Value = 0.3888021691111845
$ws.Cells.Item(6,11).Value = 3.0
$ws.Cells.Item(6,12).Value = 1.0
$ws.Cells.Item(6,13).Value = 18.76689066666667
$ws.Cells.Item(6,14).Value = 56.30067200000001
$ws.Cells.Item(6,15).Value = 0.3694347242421866
$ws.Cells.Item(6,16).Value = 0.3694347242421866
$ws.Cells.Item(6,17).Value = 28.627878299904
$ws.Cells.Item(6,18).Value = 257.650904699136
$ws.Cells.Item(6,19).Value = 0.1436370221303544
$ws.Cells.Item(6,20).Value = 0.1436370221303545

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Sema4g"
$ws.Cells.Item(7,3).Value = "Plxnb2"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3.0
$ws.Cells.Item(7,6).Value = 1.0
$ws.Cells.Item(7,7).Value = 1.525446
$ws.Cells.Item(7,8).Value = 4.576338
$ws.Cells.Item(7,9).Value = 0.3888021691111845
$ws.Cells.Item(7,10).Value = 0.3888021691111845
$ws.Cells.Item(7,11).Value = 3.0
$ws.Cells.Item(7,12).Value = 1.0
$ws.Cells.Item(7,13).Value = 23.225144
$ws.Cells.Item(7,14).Value = 69.675432
$ws.Cells.Item(7,15).Value = 0.4571974559624301
$ws.Cells.Item(7,16).Value = 0.4571974559624301
$ws.Cells.Item(7,17).Value = 35.428703014224
$ws.Cells.Item(7,18).Value = 318.858327128016
$ws.Cells.Item(7,19).Value = 0.1777593625903081
$ws.Cells.Item(7,20).Value = 0.1777593625903081

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Sema4g"
$ws.Cells.Item(8,3).Value = "Plxnb2"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3.0
$ws.Cells.Item(8,6).Value = 1.0
$ws.Cells.Item(8,7).Value = 1.96893
$ws.Cells.Item(8,8).Value = 5.90679
$ws.Cells.Item(8,9).Value = 0.5018363513543479
$ws.Cells.Item(8,10).Value = 0.5018363513543479
$ws.Cells.Item(8,11).Value = 3.0
$ws.Cells.Item(8,12).Value = 1.0
$ws.Cells.Item(8,13).Value = 8.806900666666666
$ws.Cells.Item(8,14).Value = 26.420702
$ws.Cells.Item(8,15).Value = 0.1733678197953833
$ws.Cells.Item(8,16).Value = 0.1733678197953834
$ws.Cells.Item(8,17).Value = 17.34017092962
$ws.Cells.Item(8,18).Value = 156.06153836658
$ws.Cells.Item(8,19).Value = 0.08700227412837327
$ws.Cells.Item(8,20).Value = 0.08700227412837329

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Sema4g"
$ws.Cells.Item(9,3).Value = "Plxnb2"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3.0
$ws.Cells.Item(9,6).Value = 1.0
$ws.Cells.Item(9,7).Value = 1.96893
$ws.Cells.Item(9,8).Value = 5.90679
$ws.Cells.Item(9,9).Value = 0.5018363513543479
$ws.Cells.Item(9,10).Value = 0.5018363513543479
$ws.Cells.Item(9,11).Value = 3.0
$ws.Cells.Item(9,12).Value = 1.0
$ws.Cells.Item(9,13).Value = 18.76689066666667
$ws.Cells.Item(9,14).Value = 56.30067200000001
$ws.Cells.Item(9,15).Value = 0.3694347242421866
$ws.Cells.Item(9,16).Value = 0.3694347242421866
$ws.Cells.Item(9,17).Value = 36.95069404032
$ws.Cells.Item(9,18).Value = 332.55624636288
$ws.Cells.Item(9,19).Value = 0.1853957740772986
$ws.Cells.Item(9,20).Value = 0.1853957740772986

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Sema4g"
$ws.Cells.Item(10,3).Value = "Plxnb2"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3.0
$ws.Cells.Item(10,6).Value = 1.0
$ws.Cells.Item(10,7).Value = 1.96893
$ws.Cells.Item(10,8).Value = 5.90679
$ws.Cells.Item(10,9).Value = 0.5018363513543479
$ws.Cells.Item(10,10).Value = 0.5018363513543479
$ws.Cells.Item(10,11).Value = 3.0
$ws.Cells.Item(10,12).Value = 1.0
$ws.Cells.Item(10,13).Value = 23.225144
$ws.Cells.Item(10,14).Value = 69.675432
$ws.Cells.Item(10,15).Value = 0.4571974559624301
$ws.Cells.Item(10,16).Value = 0.4571974559624301
$ws.Cells.Item(10,17).Value = 45.72868277592
$ws.Cells.Item(10,18).Value = 411.55814498328
$ws.Cells.Item(10,19).Value = 0.2294383031486761
$ws.Cells.Item(10,20).Value = 0.2294383031486761
